$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Event table v2")
$ws.Activate()

$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "User views Dashboard"
$ws.Range("C7").Value = "Data visualization"
$ws.Range("D7").Value = "System"
$ws.Range("F7").Value = "Visualized data"
$ws.Range("E7").Value = "Views Dashboard"
$ws.Range("G7").Value = "User"

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null

$ws.Range("C2:G2").Copy() | Out-Null
$ws.Range("C7:G7").PasteSpecial(-4122) | Out-Null

$ws.Range("B7").WrapText = $true
$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("B7").VerticalAlignment = -4108

$ws.Application.CutCopyMode = $false

$ws.Range("C6").Select() | Out-Null
